$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row height adjustments
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 87
$ws.Rows.Item(5).RowHeight = 87

# Fill in new row 6 data (new fine-tuned model)
$ws.Range("B6").Value = "4_bert_sc_fine_tuned"
$ws.Range("C6").Value = 8
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2e-5"
$ws.Range("D4").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("E6").Value = 3
$ws.Range("H6").Value = 0.81
$ws.Range("I6").Value = "Classified another 100 reviews from the ww2 dataset and used them tof ine-tune the model"
$ws.Range("J6").Value = "Also changed from 2 classes to three classes"

# Normalize the Notes column (I) formatting across rows 4-11 (style table cleanup)
$ws.Range("J4").Copy()
$ws.Range("I4:I11").PasteSpecial(-4122)
